$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 407
$ws.Range("I2").Value = 402.22223
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 402.22223
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -289.22223
$ws.Range("N2").Value = -676

$ws.Range("H33").Value = 6667056
$ws.Range("I33").Value = 8333524
$ws.Range("J33").Value = 1182.3334
$ws.Range("K33").Value = 8333524
$ws.Range("L33").Value = 1182.3334
$ws.Range("M33").Value = -8333295
$ws.Range("N33").Value = -1640.3334

$ws.Range("H106").Value = 2118.3333
$ws.Range("I106").Value = 1895.625
$ws.Range("K106").Value = 1895.625
$ws.Range("M106").Value = -1264.625

$ws.Range("H116").Value = 12348.4
$ws.Range("I116").Value = 14697
$ws.Range("J116").Value = 9999.799999999999
$ws.Range("K116").Value = 14697
$ws.Range("L116").Value = 9999.799999999999
$ws.Range("M116").Value = -11255
$ws.Range("N116").Value = -16883.8

$ws.Range("H133").Value = 79679.8
$ws.Range("J133").Value = 79679.8
$ws.Range("L133").Value = 79679.8
$ws.Range("N133").Value = -89799.8

$ws.Range("H137").Value = 1246.3954
$ws.Range("I137").Value = 1022.7222
$ws.Range("K137").Value = 3068.1666
$ws.Range("M137").Value = -518.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1023.6667
$ws.Range("I5").Value = 500.2
$ws.Range("J5").Value = 1397.5714
$ws.Range("K5").Value = 500.2
$ws.Range("L5").Value = 1397.5714
$ws.Range("M5").Value = -388.2
$ws.Range("N5").Value = -1621.5714

$ws.Range("H45").Value = 5889.3
$ws.Range("I45").Value = 5952.8213
$ws.Range("K45").Value = 5952.8213
$ws.Range("M45").Value = -5575.8213

$ws.Range("H74").Value = 9437.305
$ws.Range("I74").Value = 1641.5333
$ws.Range("K74").Value = 1641.5333
$ws.Range("M74").Value = -767.5333000000001

$ws.Range("H77").Value = 9437.305
$ws.Range("I77").Value = 1641.5333
$ws.Range("K77").Value = 8207.666499999999
$ws.Range("M77").Value = -3839.666499999999

$ws.Range("H97").Value = 1469.1538
$ws.Range("I97").Value = 1575
$ws.Range("J97").Value = 199
$ws.Range("K97").Value = 1575
$ws.Range("L97").Value = 199
$ws.Range("M97").Value = -1079
$ws.Range("N97").Value = -1191

$ws.Range("H102").Value = 3724.64
$ws.Range("I102").Value = 2853.9375
$ws.Range("K102").Value = 2853.9375
$ws.Range("M102").Value = -1231.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1023.6667
$ws.Range("I4").Value = 500.2
$ws.Range("J4").Value = 1397.5714
$ws.Range("K4").Value = 500.2
$ws.Range("L4").Value = 1397.5714
$ws.Range("M4").Value = -385.2
$ws.Range("N4").Value = -1627.5714

$ws.Range("H22").Value = 343.6
$ws.Range("I22").Value = 343.6
$ws.Range("K22").Value = 343.6
$ws.Range("M22").Value = -170.6

$ws.Range("H86").Value = 2914.5
$ws.Range("I86").Value = 1552.0834
$ws.Range("K86").Value = 1552.0834
$ws.Range("M86").Value = -429.0834

$ws.Range("H89").Value = 2914.5
$ws.Range("I89").Value = 1552.0834
$ws.Range("K89").Value = 7760.416999999999
$ws.Range("M89").Value = -2144.416999999999

$ws.Range("H94").Value = 1843.7368
$ws.Range("I94").Value = 1587.4857
$ws.Range("K94").Value = 1587.4857
$ws.Range("M94").Value = -1136.4857

$ws.Range("H99").Value = 5720.533
$ws.Range("I99").Value = 2224.5715
$ws.Range("J99").Value = 8779.5
$ws.Range("K99").Value = 2224.5715
$ws.Range("L99").Value = 8779.5
$ws.Range("M99").Value = -726.5715
$ws.Range("N99").Value = -11775.5

$ws.Range("H105").Value = 2307.4167
$ws.Range("I105").Value = 2187.6667
$ws.Range("K105").Value = 2187.6667
$ws.Range("M105").Value = -440.6667000000002

$ws.Range("H139").Value = 94943.25
$ws.Range("I139").Value = 80000
$ws.Range("J139").Value = 99924.336
$ws.Range("K139").Value = 80000
$ws.Range("L139").Value = 99924.336
$ws.Range("M139").Value = -74860
$ws.Range("N139").Value = -110204.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36547
$ws.Range("I31").Value = 44706.13
$ws.Range("J31").Value = 9738.429
$ws.Range("K31").Value = 44706.13
$ws.Range("L31").Value = 9738.429
$ws.Range("M31").Value = -44411.13
$ws.Range("N31").Value = -10328.429

$ws.Range("H34").Value = 36547
$ws.Range("I34").Value = 44706.13
$ws.Range("J34").Value = 9738.429
$ws.Range("K34").Value = 44706.13
$ws.Range("L34").Value = 9738.429
$ws.Range("M34").Value = -44504.13
$ws.Range("N34").Value = -10142.429

$ws.Range("H62").Value = 4149.1113
$ws.Range("J62").Value = 5973.75
$ws.Range("L62").Value = 5973.75
$ws.Range("N62").Value = -7221.75

$ws.Range("H65").Value = 4149.1113
$ws.Range("J65").Value = 5973.75
$ws.Range("L65").Value = 29868.75
$ws.Range("N65").Value = -36108.75

$ws.Range("H103").Value = 16499.5
$ws.Range("I103").Value = 10000
$ws.Range("J103").Value = 22999
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 22999
$ws.Range("M103").Value = -8828
$ws.Range("N103").Value = -25343

$ws.Range("H134").Value = 8014.9766
$ws.Range("I134").Value = 5249.5454
$ws.Range("K134").Value = 15748.6362
$ws.Range("M134").Value = -13213.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1200
$ws.Range("I59").Value = 1200
$ws.Range("K59").Value = 3600
$ws.Range("M59").Value = -3060

$ws.Range("H124").Value = 9391.357
$ws.Range("I124").Value = 7874.75
$ws.Range("K124").Value = 23624.25
$ws.Range("M124").Value = -18714.25

$ws.Range("H129").Value = 672.8333
$ws.Range("I129").Value = 672.8333
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2018.4999
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2981.5001
$ws.Range("N129").ClearContents()

$ws.Range("H131").Value = 52530.4
$ws.Range("J131").Value = 2975.5
$ws.Range("L131").Value = 8926.5
$ws.Range("N131").Value = -19006.5

$ws.Range("H134").Value = 4158.6665
$ws.Range("I134").Value = 2467.0476
$ws.Range("K134").Value = 7401.1428
$ws.Range("M134").Value = -2331.1428

$ws.Range("H137").Value = 3396
$ws.Range("I137").Value = 1693.7142
$ws.Range("J137").Value = 5779.2
$ws.Range("K137").Value = 5081.142599999999
$ws.Range("L137").Value = 17337.6
$ws.Range("M137").Value = 18.85740000000078
$ws.Range("N137").Value = -27537.6

$ws.Range("H141").Value = 207186.6
$ws.Range("J141").Value = 1000033
$ws.Range("L141").Value = 3000099
$ws.Range("N141").Value = -3010459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13713
$ws.Range("J70").Value = 15198.2
$ws.Range("L70").Value = 15198.2
$ws.Range("N70").Value = -15738.2

$ws.Range("H73").Value = 13713
$ws.Range("J73").Value = 15198.2
$ws.Range("L73").Value = 15198.2
$ws.Range("N73").Value = -17070.2

$ws.Range("H102").Value = 1654.1666
$ws.Range("I102").Value = 1654.1666
$ws.Range("K102").Value = 1654.1666
$ws.Range("M102").Value = -32.16660000000002

$ws.Range("H113").Value = 1873.0769
$ws.Range("I113").Value = 1873.0769
$ws.Range("K113").Value = 1873.0769
$ws.Range("M113").Value = 296.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3875

$ws.Range("H71").Value = 3875

$ws.Range("H93").Value = 7400
$ws.Range("J93").Value = 5000
$ws.Range("L93").Value = 5000
$ws.Range("N93").Value = -7496

$ws.Range("H132").Value = 3895.8518
$ws.Range("I132").Value = 3705.8235
$ws.Range("J132").Value = 4218.9
$ws.Range("K132").Value = 11117.4705
$ws.Range("L132").Value = 12656.7
$ws.Range("M132").Value = -8587.470499999999
$ws.Range("N132").Value = -17716.7

$ws.Range("H138").Value = 35000
$ws.Range("J138").Value = 35000
$ws.Range("L138").Value = 35000
$ws.Range("N138").Value = -45280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15512.956
$ws.Range("I62").Value = 31699.75
$ws.Range("K62").Value = 31699.75
$ws.Range("M62").Value = -31075.75

$ws.Range("H65").Value = 15512.956
$ws.Range("I65").Value = 31699.75
$ws.Range("K65").Value = 158498.75
$ws.Range("M65").Value = -155378.75

$ws.Range("H81").Value = 8261.809999999999
$ws.Range("J81").Value = 4562.533
$ws.Range("L81").Value = 9125.066000000001
$ws.Range("N81").Value = -11247.066

$ws.Range("H84").Value = 8261.809999999999
$ws.Range("J84").Value = 4562.533
$ws.Range("L84").Value = 45625.33
$ws.Range("N84").Value = -56233.33

$ws.Range("H96").Value = 1766.6666
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1766.6666
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1766.6666
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4512.6666

$ws.Range("H100").Value = 1546.8846
$ws.Range("I100").Value = 968.1429000000001
$ws.Range("K100").Value = 1936.2858
$ws.Range("M100").Value = -1395.2858

$ws.Range("H107").Value = 694.4286
$ws.Range("I107").Value = 561.06665
$ws.Range("K107").Value = 1683.19995
$ws.Range("M107").Value = 236.8000500000001

$ws.Range("H126").Value = 2643
$ws.Range("I126").Value = 2602.4167
$ws.Range("J126").Value = 2712.5715
$ws.Range("K126").Value = 7807.250100000001
$ws.Range("L126").Value = 8137.7145
$ws.Range("M126").Value = -5337.250100000001
$ws.Range("N126").Value = -13077.7145

$ws.Range("H136").Value = 2076.138
$ws.Range("I136").Value = 2132.2727
$ws.Range("J136").Value = 1899.7142
$ws.Range("K136").Value = 6396.8181
$ws.Range("L136").Value = 5699.142599999999
$ws.Range("M136").Value = -3846.8181
$ws.Range("N136").Value = -10799.1426
